# Auto-generated edit script: updates "想去人数" (F column) counts
# across sheets 展览, 演出, 本地生活, 全部类型 to match the committed snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1062
$ws.Range("F3").Value = 670
$ws.Range("F4").Value = 1479
$ws.Range("F5").Value = 3232
$ws.Range("F7").Value = 644
$ws.Range("F8").Value = 2201
$ws.Range("F9").Value = 473
$ws.Range("F10").Value = 403
$ws.Range("F13").Value = 301
$ws.Range("F14").Value = 1062
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 75
$ws.Range("F18").Value = 193
$ws.Range("F19").Value = 4386
$ws.Range("F20").Value = 1283
$ws.Range("F21").Value = 3348
$ws.Range("F22").Value = 321
$ws.Range("F23").Value = 53
$ws.Range("F24").Value = 157
$ws.Range("F25").Value = 3264
$ws.Range("F26").Value = 4870
$ws.Range("F27").Value = 122
$ws.Range("F28").Value = 969
$ws.Range("F29").Value = 538
$ws.Range("F30").Value = 3159
$ws.Range("F31").Value = 332
$ws.Range("F32").Value = 48
$ws.Range("F33").Value = 130
$ws.Range("F34").Value = 84
$ws.Range("F35").Value = 869
$ws.Range("F36").Value = 1140
$ws.Range("F37").Value = 1386
$ws.Range("F39").Value = 1311
$ws.Range("F40").Value = 836
$ws.Range("F42").Value = 785
$ws.Range("F43").Value = 490
$ws.Range("F44").Value = 49
$ws.Range("F45").Value = 276
$ws.Range("F46").Value = 55
$ws.Range("F47").Value = 138
$ws.Range("F49").Value = 3706

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 992
$ws.Range("F21").Value = 38

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2075

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2075
$ws.Range("F3").Value = 670
$ws.Range("F4").Value = 1479
$ws.Range("F5").Value = 3232
$ws.Range("F7").Value = 644
$ws.Range("F9").Value = 2201
$ws.Range("F10").Value = 473
$ws.Range("F11").Value = 403
$ws.Range("F13").Value = 992
$ws.Range("F15").Value = 301
$ws.Range("F16").Value = 1062
$ws.Range("F18").Value = 193
$ws.Range("F19").Value = 4386
$ws.Range("F21").Value = 1283
$ws.Range("F23").Value = 3348
$ws.Range("F24").Value = 3264
$ws.Range("F25").Value = 4870
$ws.Range("F26").Value = 122
$ws.Range("F27").Value = 969
$ws.Range("F28").Value = 3159
$ws.Range("F29").Value = 332
$ws.Range("F30").Value = 48
$ws.Range("F31").Value = 130
$ws.Range("F33").Value = 869
$ws.Range("F34").Value = 1140
$ws.Range("F35").Value = 1386
$ws.Range("F37").Value = 1311
$ws.Range("F39").Value = 836
$ws.Range("F40").Value = 490
$ws.Range("F42").Value = 49
$ws.Range("F44").Value = 276
$ws.Range("F46").Value = 55
$ws.Range("F47").Value = 138
$ws.Range("F49").Value = 3706
